# Auto-generated script applying scheduled market-data refresh to Kujata_Profits sheets.
# For each affected cell, write the new currentAveragePrice / LevePrice* / LeveProfit* value
# produced by the scheduled runner (see commit message).

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 449.7647
$ws.Range("I33").Value = 429.86667
$ws.Range("K33").Value = 429.86667
$ws.Range("M33").Value = -200.86667
$ws.Range("H40").Value = 1660.5
$ws.Range("J40").Value = 1367.6666
$ws.Range("L40").Value = 1367.6666
$ws.Range("N40").Value = -1717.6666
$ws.Range("H62").Value = 13893390
$ws.Range("I62").Value = 18523852
$ws.Range("K62").Value = 18523852
$ws.Range("M62").Value = -18523228
$ws.Range("H65").Value = 13893390
$ws.Range("I65").Value = 18523852
$ws.Range("K65").Value = 92619260
$ws.Range("M65").Value = -92616140
$ws.Range("H70").Value = 1750
$ws.Range("I70").Value = 1666.6666
$ws.Range("K70").Value = 4999.9998
$ws.Range("M70").Value = -4729.9998
$ws.Range("H73").Value = 1750
$ws.Range("I73").Value = 1666.6666
$ws.Range("K73").Value = 4999.9998
$ws.Range("M73").Value = -4063.9998
$ws.Range("H76").Value = 6033.1665
$ws.Range("J76").Value = 5239.8
$ws.Range("L76").Value = 5239.8
$ws.Range("N76").Value = -5869.8
$ws.Range("H79").Value = 6033.1665
$ws.Range("J79").Value = 5239.8
$ws.Range("L79").Value = 5239.8
$ws.Range("N79").Value = -7423.8
$ws.Range("H80").Value = 938.1429000000001
$ws.Range("I80").Value = 1361
$ws.Range("J80").Value = 769
$ws.Range("K80").Value = 4083
$ws.Range("L80").Value = 2307
$ws.Range("M80").Value = -3085
$ws.Range("N80").Value = -4303
$ws.Range("H83").Value = 938.1429000000001
$ws.Range("I83").Value = 1361
$ws.Range("J83").Value = 769
$ws.Range("K83").Value = 12249
$ws.Range("L83").Value = 6921
$ws.Range("M83").Value = -7257
$ws.Range("N83").Value = -16905
$ws.Range("H96").Value = 435.55554
$ws.Range("I96").Value = 284.42856
$ws.Range("J96").Value = 964.5
$ws.Range("K96").Value = 853.28568
$ws.Range("L96").Value = 2893.5
$ws.Range("M96").Value = 519.71432
$ws.Range("N96").Value = -5639.5
$ws.Range("H112").Value = 2104.8823
$ws.Range("J112").Value = 2104.8823
$ws.Range("L112").Value = 6314.646900000001
$ws.Range("N112").Value = -8530.6469
$ws.Range("H125").Value = 1883.5897
$ws.Range("I125").Value = 1707.0952
$ws.Range("J125").Value = 2089.5
$ws.Range("K125").Value = 15363.8568
$ws.Range("L125").Value = 18805.5
$ws.Range("M125").Value = -12903.8568
$ws.Range("N125").Value = -23725.5
$ws.Range("H136").Value = 47554.285
$ws.Range("J136").Value = 47554.285
$ws.Range("L136").Value = 47554.285
$ws.Range("N136").Value = -57754.285
$ws.Range("H137").Value = 1093.875
$ws.Range("J137").Value = 1871.1
$ws.Range("L137").Value = 5613.299999999999
$ws.Range("N137").Value = -10713.3
$ws.Range("H138").Value = 2834.443
$ws.Range("I138").Value = 2457.1667
$ws.Range("J138").Value = 2902.015
$ws.Range("K138").Value = 7371.500100000001
$ws.Range("L138").Value = 8706.045
$ws.Range("M138").Value = -2231.500100000001
$ws.Range("N138").Value = -18986.045

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9939.139999999999
$ws.Range("I32").Value = 7971.355
$ws.Range("J32").Value = 13149.737
$ws.Range("K32").Value = 7971.355
$ws.Range("L32").Value = 13149.737
$ws.Range("M32").Value = -7684.355
$ws.Range("N32").Value = -13723.737
$ws.Range("H74").Value = 1173.26
$ws.Range("I74").Value = 603
$ws.Range("J74").Value = 2280.2354
$ws.Range("K74").Value = 603
$ws.Range("L74").Value = 2280.2354
$ws.Range("M74").Value = 271
$ws.Range("N74").Value = -4028.2354
$ws.Range("H77").Value = 1173.26
$ws.Range("I77").Value = 603
$ws.Range("J77").Value = 2280.2354
$ws.Range("K77").Value = 3015
$ws.Range("L77").Value = 11401.177
$ws.Range("M77").Value = 1353
$ws.Range("N77").Value = -20137.177
$ws.Range("H97").Value = 6395.5293
$ws.Range("I97").Value = 552.25
$ws.Range("K97").Value = 552.25
$ws.Range("M97").Value = -56.25
$ws.Range("H102").Value = 7939222
$ws.Range("I102").Value = 7939222
$ws.Range("K102").Value = 7939222
$ws.Range("M102").Value = -7937600
$ws.Range("H132").Value = 2317.093
$ws.Range("I132").Value = 1878.2
$ws.Range("K132").Value = 5634.6
$ws.Range("M132").Value = -3104.6

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 45455696
$ws.Range("I99").Value = 58824576
$ws.Range("J99").Value = 1510
$ws.Range("K99").Value = 58824576
$ws.Range("L99").Value = 1510
$ws.Range("M99").Value = -58823078
$ws.Range("N99").Value = -4506
$ws.Range("H105").Value = 55557036
$ws.Range("I105").Value = 58824920
$ws.Range("J105").Value = 3000
$ws.Range("K105").Value = 58824920
$ws.Range("L105").Value = 3000
$ws.Range("M105").Value = -58823173
$ws.Range("N105").Value = -6494
$ws.Range("H107").Value = 892.8570999999999
$ws.Range("I107").Value = 902.5
$ws.Range("J107").Value = 700
$ws.Range("K107").Value = 902.5
$ws.Range("L107").Value = 700
$ws.Range("M107").Value = 1017.5
$ws.Range("N107").Value = -4540
$ws.Range("H134").Value = 4543.3667
$ws.Range("I134").Value = 978.13336
$ws.Range("K134").Value = 2934.40008
$ws.Range("M134").Value = -399.4000800000003

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1708.7222
$ws.Range("I31").Value = 1665.44
$ws.Range("J31").Value = 2249.75
$ws.Range("K31").Value = 1665.44
$ws.Range("L31").Value = 2249.75
$ws.Range("M31").Value = -1370.44
$ws.Range("N31").Value = -2839.75
$ws.Range("H34").Value = 1708.7222
$ws.Range("I34").Value = 1665.44
$ws.Range("J34").Value = 2249.75
$ws.Range("K34").Value = 1665.44
$ws.Range("L34").Value = 2249.75
$ws.Range("M34").Value = -1463.44
$ws.Range("N34").Value = -2653.75
$ws.Range("H109").Value = 17633.666
$ws.Range("J109").Value = 17633.666
$ws.Range("L109").Value = 17633.666
$ws.Range("N109").Value = -19713.666
$ws.Range("H133").Value = 39960
$ws.Range("J133").Value = 39960
$ws.Range("L133").Value = 39960
$ws.Range("N133").Value = -45020
$ws.Range("H134").Value = 8621969
$ws.Range("I134").Value = 1179.3903
$ws.Range("J134").Value = 29413284
$ws.Range("K134").Value = 3538.1709
$ws.Range("L134").Value = 88239852
$ws.Range("M134").Value = -1003.1709
$ws.Range("N134").Value = -88244922
$ws.Range("H141").Value = 308359.8
$ws.Range("J141").Value = 308359.8
$ws.Range("L141").Value = 308359.8
$ws.Range("N141").Value = -318719.8

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 3402.182
$ws.Range("I81").Value = 2399.75
$ws.Range("J81").Value = 3975
$ws.Range("K81").Value = 7199.25
$ws.Range("L81").Value = 11925
$ws.Range("M81").Value = -6076.25
$ws.Range("N81").Value = -14171
$ws.Range("H84").Value = 3402.182
$ws.Range("I84").Value = 2399.75
$ws.Range("J84").Value = 3975
$ws.Range("K84").Value = 21597.75
$ws.Range("L84").Value = 35775
$ws.Range("M84").Value = -15981.75
$ws.Range("N84").Value = -47007
$ws.Range("H131").Value = 28615234
$ws.Range("I131").Value = 125000410
$ws.Range("J131").Value = 56664.965
$ws.Range("K131").Value = 375001230
$ws.Range("L131").Value = 169994.895
$ws.Range("M131").Value = -374996190
$ws.Range("N131").Value = -180074.895

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H105").Value = 20000
$ws.Range("J105").Value = 20000
$ws.Range("L105").Value = 20000
$ws.Range("N105").Value = -26988
$ws.Range("H126").Value = 1948
$ws.Range("I126").Value = 1320
$ws.Range("K126").Value = 3960
$ws.Range("M126").Value = -1490

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3335
$ws.Range("I40").Value = 3335
$ws.Range("K40").Value = 3335
$ws.Range("M40").Value = -3199
$ws.Range("H46").Value = 1790.4
$ws.Range("J46").Value = 1790.4
$ws.Range("L46").Value = 1790.4
$ws.Range("N46").Value = -2166.4
$ws.Range("H55").Value = 372.58334
$ws.Range("J55").Value = 544.8
$ws.Range("L55").Value = 544.8
$ws.Range("N55").Value = -890.8
$ws.Range("H68").Value = 956.25
$ws.Range("I68").Value = 826
$ws.Range("J68").Value = 1086.5
$ws.Range("K68").Value = 826
$ws.Range("L68").Value = 1086.5
$ws.Range("M68").Value = -77
$ws.Range("N68").Value = -2584.5
$ws.Range("H71").Value = 956.25
$ws.Range("I71").Value = 826
$ws.Range("J71").Value = 1086.5
$ws.Range("K71").Value = 4130
$ws.Range("L71").Value = 5432.5
$ws.Range("M71").Value = -386
$ws.Range("N71").Value = -12920.5
$ws.Range("H93").Value = 1008.6429
$ws.Range("I93").Value = 961.0909
$ws.Range("K93").Value = 961.0909
$ws.Range("M93").Value = 286.9091
$ws.Range("H132").Value = 3069.5715
$ws.Range("I132").Value = 2568.1428
$ws.Range("J132").Value = 3571
$ws.Range("K132").Value = 7704.428400000001
$ws.Range("L132").Value = 10713
$ws.Range("M132").Value = -5174.428400000001
$ws.Range("N132").Value = -15773
$ws.Range("H136").Value = 2301.125
$ws.Range("I136").Value = 1201.3334
$ws.Range("J136").Value = 2961
$ws.Range("K136").Value = 3604.0002
$ws.Range("L136").Value = 8883
$ws.Range("M136").Value = -1054.0002
$ws.Range("N136").Value = -13983

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 50342
$ws.Range("I109").Value = 50342
$ws.Range("K109").Value = 50342
$ws.Range("M109").Value = -48955
$ws.Range("H139").Value = 38132
$ws.Range("J139").Value = 38132
$ws.Range("L139").Value = 38132
$ws.Range("N139").Value = -48412
